$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (report volume/number + week-covering date range)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  4/29/2024  Through  5/5/2024"

# ---------------------------------------------------------------------------
# Crime-complaints table (rows 15-31), columns:
#   C = Week to Date 2024, D = Week to Date 2023, E = %Chg
#   F = 28 Day 2024,       G = 28 Day 2023,        H = %Chg
#   I = Year to Date 2024, J = Year to Date 2023,  K = %Chg
#   L = 2 Year %Chg, M = 14 Year %Chg, N = 31 Year %Chg
# ---------------------------------------------------------------------------

# Row 15: Rape -- C15 & F15 flip from the text placeholder "0" to a real number
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 1
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("I15").Value = 4
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = 100
$ws.Range("M15").Value = 33.333333333333
$ws.Range("N15").Value = -63.636363636363

# Row 16: Robbery
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -60
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -20
$ws.Range("I16").Value = 65
$ws.Range("J16").Value = 58
$ws.Range("K16").Value = 12.068965517241
$ws.Range("L16").Value = 10.169491525423
$ws.Range("M16").Value = -5.797101449275
$ws.Range("N16").Value = -75.746268656716

# Row 17: Fel. Assault
$ws.Range("C17").Value = 6
$ws.Range("E17").Value = 200
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 46.666666666666
$ws.Range("I17").Value = 88
$ws.Range("J17").Value = 87
$ws.Range("K17").Value = 1.149425287356
$ws.Range("L17").Value = 8.641975308641
$ws.Range("M17").Value = 39.682539682539
$ws.Range("N17").Value = -57.073170731707

# Row 18: Burglary
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -62.5
$ws.Range("I18").Value = 31
$ws.Range("J18").Value = 43
$ws.Range("K18").Value = -27.906976744186
$ws.Range("L18").Value = -53.731343283582
$ws.Range("M18").Value = -24.390243902439
$ws.Range("N18").Value = -89.836065573770

# Row 19: Gr. Larceny
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 266.666666666667
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = 164.705882352941
$ws.Range("I19").Value = 177
$ws.Range("J19").Value = 118
$ws.Range("K19").Value = 50
$ws.Range("L19").Value = 48.739495798319
$ws.Range("M19").Value = 82.474226804123
$ws.Range("N19").Value = 25.531914893617

# Row 20: G.L.A.
$ws.Range("D20").Value = 4
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -80
$ws.Range("J20").Value = 20
$ws.Range("K20").Value = -25
$ws.Range("L20").Value = -6.25
$ws.Range("N20").Value = -72.727272727272

# Row 21: TOTAL
$ws.Range("C21").Value = 21
$ws.Range("E21").Value = 40
$ws.Range("F21").Value = 84
$ws.Range("G21").Value = 61
$ws.Range("H21").Value = 37.704918032786
$ws.Range("I21").Value = 382
$ws.Range("J21").Value = 331
$ws.Range("K21").Value = 15.407854984894
$ws.Range("L21").Value = 10.724637681159
$ws.Range("M21").Value = 35.943060498220
$ws.Range("N21").Value = -61.685055165496

# Row 22: Transit
$ws.Range("M22").Value = -25

# Row 23: Housing
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 2
$ws.Range("F23").Value = 7
$ws.Range("H23").Value = 40
$ws.Range("I23").Value = 33
$ws.Range("J23").Value = 31
$ws.Range("K23").Value = 6.451612903225
$ws.Range("L23").Value = -2.941176470588
$ws.Range("M23").Value = 200

# Row 24: Petit Larceny
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 66.666666666666
$ws.Range("F24").Value = 143
$ws.Range("G24").Value = 100
$ws.Range("H24").Value = 43
$ws.Range("I24").Value = 433
$ws.Range("J24").Value = 407
$ws.Range("K24").Value = 6.388206388206
$ws.Range("L24").Value = 16.711590296496
$ws.Range("M24").Value = 39.677419354838

# Row 25: Retail Theft
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 27.272727272727
$ws.Range("F25").Value = 101
$ws.Range("G25").Value = 54
$ws.Range("H25").Value = 87.037037037037
$ws.Range("I25").Value = 276
$ws.Range("J25").Value = 251
$ws.Range("K25").Value = 9.960159362549
$ws.Range("L25").Value = 41.538461538461

# Row 26: Misd. Assault
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 36
$ws.Range("H26").Value = -16.666666666666
$ws.Range("I26").Value = 128
$ws.Range("J26").Value = 137
$ws.Range("K26").Value = -6.569343065693
$ws.Range("L26").Value = -17.948717948717
$ws.Range("M26").Value = -23.809523809523

# Row 27: UCR Rape* -- C27/D27/F27 flip from the text placeholder "0" to real numbers,
#         E27 flips from the text placeholder "***.*" to a real percentage number
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = 0
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F27").Value = 1
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 6
$ws.Range("J27").Value = 5
$ws.Range("K27").Value = 20
$ws.Range("L27").Value = 50

# Row 28: Other Sex Crimes
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 2
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 16
$ws.Range("J28").Value = 14
$ws.Range("K28").Value = 14.285714285714
$ws.Range("L28").Value = -23.809523809523

# Row 29: Shooting Vic.
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("N29").Value = -81.818181818181

# Row 30: Shooting Inc.
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("N30").Value = -82.758620689655

# Row 31: Hate Crimes -- G31/H31 flip from real numbers back to the text placeholders
#         "0" / "***.*" (no incidents recorded this period). Copy the formatting+content
#         from existing cells that already carry that exact placeholder style/text so no
#         stray number-format gets registered in the workbook's style table.
$ws.Range("C14").Copy($ws.Range("G31"))
$ws.Range("E15").Copy($ws.Range("H31"))
